# Save-the-Salmon-Models / lostine.xlsx -> Chinook-Life-Cycle edit
# "updated model, added weir fun, and included draft write-up methods"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Updated model numbers (cols M-Q on rows 17,18,20-25) -----------------
$ws.Range("M17").Value = 181
$ws.Range("N17").Value = 1271

$ws.Range("M18").Value = 87

$ws.Range("Q20").Value = 519

$ws.Range("N21").Value = 285
$ws.Range("P21").Value = 177
$ws.Range("Q21").Value = 726

$ws.Range("M22").Value = 13
$ws.Range("N22").Value = 121
$ws.Range("P22").Value = 62

$ws.Range("M23").Value = 45
$ws.Range("N23").Value = 296
$ws.Range("P23").Value = 140
$ws.Range("Q23").Value = 591

$ws.Range("M24").Value = 7
$ws.Range("P24").Value = 209
$ws.Range("Q24").Value = 586

$ws.Range("P25").Value = 137

# --- Column widths: split the uniform M:R band into individually sized ---
# --- columns (weir-count / fun-count columns got their own widths) -------
$ws.Columns("M").ColumnWidth = 12.166666666666666
$ws.Columns("N").ColumnWidth = 9.166666666666666
$ws.Columns("O").ColumnWidth = 9.022135416666666
$ws.Columns("P").ColumnWidth = 9.592447916666666
$ws.Columns("Q").ColumnWidth = 7.022135416666667
$ws.Columns("R").ColumnWidth = 7.307291666666667

# --- Selection moved from L11 to B1 --------------------------------------
[void]$ws.Range("B1").Select()
